$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVTs")
$ws.Activate()

# Copy formatting from the row above (row 7) so the new rows match the
# existing table styling (borders, fonts, wrap text, etc.)
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# E8/E9 use the left-aligned wrap style (same as E5) instead of the default E-column style.
$ws.Range("E5").Copy()
$ws.Range("E8:E9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 8: Legend Interactivity ---
$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = "Legend Interactivity"
$ws.Cells.Item(8, 4).Value = "Click on Legend value."
$ws.Cells.Item(8, 5).Value = "Data should get filtered on clicking particular legend value. Also opacity should get changed for the circle accordingly."
$ws.Cells.Item(8, 3).Value = "Check whether Legend is interactive."

# --- Row 9: Visual Interactivity ---
$ws.Cells.Item(9, 2).Value = "Visual Interactivity"
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 3).Value = "Check whether Visual is interactive."
$ws.Cells.Item(9, 5).Value = "Data should get filtered for the ""Maths"" category. Also opacity should get changed for that circle accordingly."
$ws.Cells.Item(9, 4).Value = "Click on ""Maths""  path element(circle).                          [NOTE : For overlapping(Intersection) path elements there is no interactivity]                               "

# Row heights
$ws.Rows.Item(8).RowHeight = 135.75
$ws.Rows.Item(9).RowHeight = 105

# Update dimension / selection to match new extent
$ws.Range("D9").Select()
